# Updating the base model for Astro
# - Shift the "Data" (date/time) column A from 20.06.2024 -> 25.06.2024 (+5 days)
# - Shift the "Lookup" column D text from "20.06.2024<n>" -> "25.06.2024<n>"
# - Update several "Prediction" values in column C (rows 27-86)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Prediction (column C) values for rows 27-86 (1-indexed worksheet rows)
$newC = @{
    27 = 0.015
    28 = 0.03
    29 = 0.045
    30 = 0.068
    31 = 0.096
    32 = 0.142
    33 = 0.187
    34 = 0.22
    35 = 0.274
    36 = 0.332
    37 = 0.378
    38 = 0.418
    39 = 0.46
    40 = 0.506
    41 = 0.553
    42 = 0.5649999999999999
    43 = 0.605
    44 = 0.671
    45 = 0.704
    46 = 0.734
    47 = 0.753
    48 = 0.787
    49 = 0.79
    50 = 0.799
    51 = 0.822
    52 = 0.843
    53 = 0.843
    54 = 0.843
    55 = 0.843
    56 = 0.843
    57 = 0.831
    58 = 0.8169999999999999
    59 = 0.801
    60 = 0.771
    61 = 0.6879999999999999
    62 = 0.605
    63 = 0.59
    64 = 0.553
    65 = 0.517
    66 = 0.443
    67 = 0.415
    68 = 0.395
    69 = 0.394
    70 = 0.315
    71 = 0.292
    72 = 0.251
    73 = 0.211
    74 = 0.176
    75 = 0.148
    76 = 0.113
    77 = 0.094
    78 = 0.076
    79 = 0.077
    80 = 0.063
    81 = 0.055
    82 = 0.049
    83 = 0.037
    84 = 0.027
    85 = 0.019
    86 = 0.014
}

for ($row = 2; $row -le 96; $row++) {
    # Column A: shift the date/time value forward by exactly 5 days (serial date + 5)
    $oldDate = $ws.Cells.Item($row, 1).Value2
    $ws.Cells.Item($row, 1).Value2 = $oldDate + 5

    # Column D: rewrite the lookup text, replacing the old date prefix with the new one
    $oldLookup = $ws.Cells.Item($row, 4).Value2
    $newLookup = $oldLookup -replace '^20\.06\.2024', '25.06.2024'
    $ws.Cells.Item($row, 4).Value2 = $newLookup

    # Column C: apply updated prediction values where changed
    if ($newC.ContainsKey($row)) {
        $ws.Cells.Item($row, 3).Value2 = $newC[$row]
    }
}
